# Update stats for 2026-02 (row 27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6548
$ws.Range("D27").Value = 6105225
$ws.Range("E27").Value = 932.3801160659743
$ws.Range("F27").Value = 10.05042016806723
$ws.Range("H27").Value = 25.22236148978967
